$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Identity" column (currently column D).
# Excel's default column-insert behavior copies formatting from the column to
# the left (column C), which matches the target style pattern exactly.
$ws.Columns("D:D").Insert()

# Populate the new "format" column. Write order matches the original
# authoring session so new shared-string entries land in the same order
# (alphabets, normal, format).
$ws.Range("D3").Value = "alphabets"
$ws.Range("D4").Value = "alphabets"
$ws.Range("D5").Value = "alphabets"
$ws.Range("D2").Value = "normal"
$ws.Range("D1").Value = "format"

# Adjust column widths: column B (now narrower-fit) and the new column D.
$ws.Columns("B:B").ColumnWidth = 5.5
$ws.Columns("D:D").ColumnWidth = 12.1

# Update the selection to match the saved view state.
$ws.Range("D4").Select()
